$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The display_id for the 2nd data row (originally "cs0002_slr0612") is now
# constructed from the key at transform-time, so the literal value becomes
# the templated "cs0002_{key}" and the key itself ("slr0612") is populated
# into the new key column for that row.
$ws.Range("C3").Value = "slr0612"
$ws.Range("A3").Value = "cs0002_{key}"

$ws.Range("H10").Select() | Out-Null
